# Split the single run
#   "— Hit Ctrl+P (or Cmd+P on Mac) on any page for quick copies without downloading"
# into five runs:
#   "— Hit" / " " / "Ctrl+P" / " " / "(or Cmd+P on Mac) on any page for quick copies without downloading"
# without altering the visible text or formatting.
#
# The run-splitting is done by adding a (zero-content-impact) comment on each
# internal boundary sub-range and immediately deleting it again: Word always
# gives a commented sub-range its own run, and — unlike toggling a character
# property such as Bold on/off — removing the comment afterwards leaves no
# stray <w:rPr/> behind on the resulting runs.

$d = $word.ActiveDocument

$full = "— Hit Ctrl+P (or Cmd+P on Mac) on any page for quick copies without downloading"

$r = $d.Content
$found = $r.Find.Execute($full, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target 'Print from browser' sentence."
}

$start = $r.Start
$end = $r.End

# Internal split boundaries (relative to $start), matching:
#   0   "— Hit"
#   5   " "
#   6   "Ctrl+P"
#   12  " "
#   13  "(or Cmd+P on Mac) on any page for quick copies without downloading"
#   79  (end)
$boundaryRanges = @(
    @($start + 5, $start + 6),
    @($start + 6, $start + 12),
    @($start + 12, $start + 13)
)

foreach ($b in $boundaryRanges) {
    $sub = $d.Range($b[0], $b[1])
    $cmt = $d.Comments.Add($sub, "split")
    $cmt.Delete()
}
